{"js": "// Update the date title and the 25 division-problem answers in the table.\n// Mapping is strictly positional (title paragraph, then table rows/cols in\n// reading order) so it is immune to duplicate text values (e.g. \"21\u00f79=2, 3\"\n// occurs twice in the output).\n\nconst titleText = \"2024-03-11 Monday\";\n\n// Row-major list of the 25 new cell values (5 rows x 5 cols actually used;\n// the table has 20 rows total but only rows 0, 4, 8, 12, 16 contain text).\nconst newCellValues = [\n  [\"77\u00f74=19, 1\", \"72\u00f72=36, 0\", \"71\u00f76=11, 5\", \"44\u00f78=5, 4\", \"84\u00f72=42, 0\"],\n  [\"40\u00f79=4, 4\", \"73\u00f77=10, 3\", \"65\u00f73=21, 2\", \"21\u00f79=2, 3\", \"16\u00f78=2, 0\"],\n  [\"14\u00f77=2, 0\", \"33\u00f76=5, 3\", \"42\u00f77=6, 0\", \"17\u00f72=8, 1\", \"21\u00f79=2, 3\"],\n  [\"14\u00f76=2, 2\", \"19\u00f78=2, 3\", \"92\u00f78=11, 4\", \"54\u00f75=10, 4\", \"27\u00f78=3, 3\"],\n  [\"28\u00f77=4, 0\", \"34\u00f77=4, 6\", \"43\u00f78=5, 3\", \"35\u00f78=4, 3\", \"35\u00f75=7, 0\"],\n];\n\n// The physical rows (of the 20-row table) that actually hold the answers.\nconst dataRowIndexes = [0, 4, 8, 12, 16];\n\nconst body = context.document.body;\n\n// --- Update the title paragraph (first paragraph in the body) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(titleText, Word.InsertLocation.replace);\n\n// --- Update the table cells ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < dataRowIndexes.length; r++) {\n  const physicalRow = dataRowIndexes[r];\n  for (let c = 0; c < newCellValues[r].length; c++) {\n    const cell = table.getCell(physicalRow, c);\n    cell.value = newCellValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and the 25 division-problem answers in the table.\n# Cell coordinates are used (1-indexed, per Word COM convention) instead of\n# Find/Replace so that the duplicate new value \"21\u00f79=2, 3\" (which appears\n# twice in the result) cannot cause an incorrect/ambiguous replacement.\n\n$d = $word.ActiveDocument\n\n# --- Update the title paragraph (first paragraph in the document) ---\n$d.Paragraphs.Item(1).Range.Text = \"2024-03-11 Monday\"\n\n# --- Update the table cells ---\n$table = $d.Tables.Item(1)\n\n# Physical rows (1-indexed) of the 20-row table that hold answers, and the\n# new value for each of the 5 columns in that row, in order.\n$newValues = @{\n    1  = @(\"77\u00f74=19, 1\", \"72\u00f72=36, 0\", \"71\u00f76=11, 5\", \"44\u00f78=5, 4\", \"84\u00f72=42, 0\")\n    5  = @(\"40\u00f79=4, 4\", \"73\u00f77=10, 3\", \"65\u00f73=21, 2\", \"21\u00f79=2, 3\", \"16\u00f78=2, 0\")\n    9  = @(\"14\u00f77=2, 0\", \"33\u00f76=5, 3\", \"42\u00f77=6, 0\", \"17\u00f72=8, 1\", \"21\u00f79=2, 3\")\n    13 = @(\"14\u00f76=2, 2\", \"19\u00f78=2, 3\", \"92\u00f78=11, 4\", \"54\u00f75=10, 4\", \"27\u00f78=3, 3\")\n    17 = @(\"28\u00f77=4, 0\", \"34\u00f77=4, 6\", \"43\u00f78=5, 3\", \"35\u00f78=4, 3\", \"35\u00f75=7, 0\")\n}\n\nforeach ($row in $newValues.Keys) {\n    $cols = $newValues[$row]\n    for ($c = 0; $c -lt $cols.Length; $c++) {\n        $table.Cell($row, $c + 1).Range.Text = $cols[$c]\n    }\n}\n"}
